{"js": "// Replace the two-digit-multiplication \"fact\" strings in the table with\n// their updated answers. Each old string is unique in the document, so we\n// locate it with a body search and replace the whole run of text in one\n// shot via insertText(..., \"Replace\").\n//\n// NOTE on ordering: one new value (\"83\u00d712=996\") happens to collide with an\n// old value that appears elsewhere in the table (and that old value is\n// itself later replaced by a different new value, \"47\u00d799=4653\"). If we\n// replaced in naive top-to-bottom order, the earlier replacement would\n// manufacture a second \"83\u00d712=996\" in the document, and the later search\n// for \"83\u00d712=996\" would then match two cells instead of one. To avoid\n// that, the pair that originally reads \"83\u00d712=996\" is replaced BEFORE the\n// pair that is replaced WITH \"83\u00d712=996\".\nconst replacements = [\n  [\"63\u00d731=1953\", \"58\u00d761=3538\"],\n  [\"83\u00d729=2407\", \"60\u00d739=2340\"],\n  [\"73\u00d797=7081\", \"46\u00d781=3726\"],\n  [\"83\u00d731=2573\", \"31\u00d736=1116\"],\n  [\"31\u00d793=2883\", \"48\u00d774=3552\"],\n  [\"47\u00d793=4371\", \"21\u00d794=1974\"],\n  [\"61\u00d717=1037\", \"77\u00d713=1001\"],\n  [\"87\u00d739=3393\", \"22\u00d758=1276\"],\n  [\"80\u00d776=6080\", \"50\u00d746=2300\"],\n  [\"12\u00d794=1128\", \"38\u00d741=1558\"],\n  [\"84\u00d729=2436\", \"88\u00d743=3784\"],\n  [\"80\u00d793=7440\", \"35\u00d790=3150\"],\n  [\"83\u00d712=996\", \"47\u00d799=4653\"],   // must run before \"39\u00d724=936\" -> \"83\u00d712=996\"\n  [\"39\u00d724=936\", \"83\u00d712=996\"],\n  [\"66\u00d719=1254\", \"61\u00d751=3111\"],\n  [\"72\u00d720=1440\", \"78\u00d783=6474\"],\n  [\"35\u00d741=1435\", \"32\u00d749=1568\"],\n  [\"81\u00d714=1134\", \"78\u00d732=2496\"],\n  [\"28\u00d774=2072\", \"66\u00d741=2706\"],\n  [\"55\u00d759=3245\", \"70\u00d788=6160\"],\n  [\"19\u00d775=1425\", \"25\u00d744=1100\"],\n  [\"26\u00d761=1586\", \"78\u00d784=6552\"],\n  [\"78\u00d774=5772\", \"31\u00d789=2759\"],\n  [\"77\u00d781=6237\", \"68\u00d756=3808\"],\n  [\"71\u00d799=7029\", \"39\u00d766=2574\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match for \"${oldText}\", found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit-multiplication \"fact\" strings in the table with\n# their updated answers.\n#\n# NOTE on ordering: one new value (\"83\u00d712=996\") happens to collide with an\n# old value that appears elsewhere in the table (and that old value is\n# itself later replaced by a different new value, \"47\u00d799=4653\"). If we\n# replaced in naive top-to-bottom order, the earlier replacement would\n# manufacture a second \"83\u00d712=996\" in the document, and the later\n# find/replace for \"83\u00d712=996\" would then match (and stomp) two cells\n# instead of one. To avoid that, the pair that originally reads\n# \"83\u00d712=996\" is replaced BEFORE the pair that is replaced WITH\n# \"83\u00d712=996\" below.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"63\u00d731=1953\", \"58\u00d761=3538\")\n    ,@(\"83\u00d729=2407\", \"60\u00d739=2340\")\n    ,@(\"73\u00d797=7081\", \"46\u00d781=3726\")\n    ,@(\"83\u00d731=2573\", \"31\u00d736=1116\")\n    ,@(\"31\u00d793=2883\", \"48\u00d774=3552\")\n    ,@(\"47\u00d793=4371\", \"21\u00d794=1974\")\n    ,@(\"61\u00d717=1037\", \"77\u00d713=1001\")\n    ,@(\"87\u00d739=3393\", \"22\u00d758=1276\")\n    ,@(\"80\u00d776=6080\", \"50\u00d746=2300\")\n    ,@(\"12\u00d794=1128\", \"38\u00d741=1558\")\n    ,@(\"84\u00d729=2436\", \"88\u00d743=3784\")\n    ,@(\"80\u00d793=7440\", \"35\u00d790=3150\")\n    ,@(\"83\u00d712=996\",  \"47\u00d799=4653\")   # must run before \"39\u00d724=936\" -> \"83\u00d712=996\"\n    ,@(\"39\u00d724=936\",  \"83\u00d712=996\")\n    ,@(\"66\u00d719=1254\", \"61\u00d751=3111\")\n    ,@(\"72\u00d720=1440\", \"78\u00d783=6474\")\n    ,@(\"35\u00d741=1435\", \"32\u00d749=1568\")\n    ,@(\"81\u00d714=1134\", \"78\u00d732=2496\")\n    ,@(\"28\u00d774=2072\", \"66\u00d741=2706\")\n    ,@(\"55\u00d759=3245\", \"70\u00d788=6160\")\n    ,@(\"19\u00d775=1425\", \"25\u00d744=1100\")\n    ,@(\"26\u00d761=1586\", \"78\u00d784=6552\")\n    ,@(\"78\u00d774=5772\", \"31\u00d789=2759\")\n    ,@(\"77\u00d781=6237\", \"68\u00d756=3808\")\n    ,@(\"71\u00d799=7029\", \"39\u00d766=2574\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1          # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $true\n\n    $ok = $find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $ok) {\n        Write-Output \"WARNING: replacement not found for '$oldText'\"\n    }\n}\n\nWrite-Output $d.Content.Text\n"}
